# Fix the Abstract (D5) and Authors (E5) cells for row 5 of the references
# sheet: the abstract text had stray id="ParN"> fragments and doubled blank
# lines left over from the source scrape, and the authors list needs one
# more level of the repeated-whitespace-separator bug applied (matching the
# pattern already present across the other duplicated author/abstract rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$abstract = @"
Background
Recently, the World Health Organization has declared the coronavirus disease 2019 (COVID-19) outbreak a public health emergency of international concern.
 So far, however, limited data are available for children.
 Therefore, we aimed to investigate the clinical and chest CT imaging characteristics of COVID-19 in preschool children.
Methods
From January 26, 2020 to February 20, 2020, the clinical and initial chest CT imaging data of eight preschool children with laboratory-confirmed COVID-19 from two hospitals were retrospectively collected.
 The chest CT imaging characteristics, including the distribution, shape, and density of lesions, and the pleural effusion, pleural changes, and enlarged lymph nodes were evaluated.
Results
Two cases (25%) were classified as mild type, and they showed no obvious abnormal CT findings or minimal pleural thickening on the right side.
 Five cases (62.5%) were classified as moderate type.
 Among these patients, one case showed consolidation located in the subpleural region of the right upper lobe, with thickening in the adjacent pleura; one case showed multiple consolidation and ground-glass opacities with blurry margins; one case displayed bronchial pneumonia-like changes in the left upper lobe; and two cases displayed asthmatic bronchitis-like changes.
 One case (12.5%) was classified as critical type and showed bronchial pneumonia-like changes in the bilateral lungs, presenting blurred and messy bilateral lung markings and multiple patchy shadows scattered along the lung markings with blurry margins.
Conclusions
The chest CT findings of COVID-19 in preschool children are atypical and various.
 Accurate diagnosis requires a comprehensive evaluation of epidemiological, clinical, laboratory and CT imaging data.

"@

$authors = @"
[Yang%Li%NULL%0,        Jianghui%Cao%NULL%1,        Xiaolong%Zhang%NULL%1,        Guangzhi%Liu%NULL%1,        Xiaxia%Wu%NULL%1,        Baolin%Wu%WBLlin0129@163.com%1]
"@

$ws.Range("D5").Value = $abstract
$ws.Range("E5").Value = $authors
